# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect newly scraped counts (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2147
$ws1.Range("F4").Value = 29
$ws1.Range("F5").Value = 11240
$ws1.Range("F6").Value = 197
$ws1.Range("F7").Value = 172
$ws1.Range("F10").Value = 11162
$ws1.Range("F13").Value = 53
$ws1.Range("F14").Value = 1730
$ws1.Range("F15").Value = 5588
$ws1.Range("F16").Value = 97

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 2147
$ws4.Range("F5").Value = 29
$ws4.Range("F7").Value = 11240
$ws4.Range("F8").Value = 197
$ws4.Range("F9").Value = 172
$ws4.Range("F12").Value = 11162
$ws4.Range("F15").Value = 53
$ws4.Range("F16").Value = 1730
$ws4.Range("F17").Value = 5588
$ws4.Range("F18").Value = 97
